$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Title: "Metrics Summary" -> "Executive Summary" ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Executive Summary"

# --- Content placeholder: replace the long markdown dump with a concise summary ---
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "`rTotal GLA: 222,221 m²`rOccupancy: 100% (based on total leased area matching GLA)`rWALT: Not directly computable from provided data (requires weighted average calculation of lease terms)`rIn-Place Rent: £5.5 per m² per annum (based on the lease with Ingram Micro)`rKey Highlight 1: Asset is a logistics facility with significant parking capacity (180 spaces)`rKey Highlight 2: Located in Daventry, UK, with a strategic logistics position"
